$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.484
$ws.Range("B7").Value = 5.525
$ws.Range("A10").Value = -21.252
$ws.Range("D10").Value = -7.926
$ws.Range("A12").Value = -21.515
$ws.Range("D14").Value = -7.637
$ws.Range("B15").Value = 5.236000000000001
$ws.Range("E16").Value = 16.747
$ws.Range("A18").Value = -21.866
$ws.Range("C18").Value = -11.081
$ws.Range("C19").Value = -12.051
$ws.Range("B20").Value = 7.217000000000001
$ws.Range("E22").Value = 16.581
$ws.Range("E26").Value = 16.753
$ws.Range("C27").Value = -13.301
$ws.Range("B29").Value = 5.626
$ws.Range("B30").Value = 5.131
$ws.Range("B31").Value = 6.114
$ws.Range("D32").Value = -8.274000000000001
$ws.Range("D35").Value = -7.387
$ws.Range("A37").Value = -20.02
$ws.Range("B40").Value = 9.208
$ws.Range("C42").Value = -12.328
$ws.Range("D43").Value = -8.375
$ws.Range("C44").Value = -12.536
$ws.Range("E44").Value = 16.612
$ws.Range("C47").Value = -12.24
$ws.Range("D49").Value = -7.997000000000002
$ws.Range("E54").Value = 16.686
$ws.Range("A55").Value = -21.875
$ws.Range("D56").Value = -8.119
$ws.Range("C58").Value = -12.575
$ws.Range("E63").Value = 17.565
$ws.Range("A68").Value = -21.581
$ws.Range("B68").Value = 5.242
$ws.Range("D69").Value = -6.812
$ws.Range("C73").Value = -12.77
$ws.Range("B76").Value = 6.175
$ws.Range("A77").Value = -20.252
$ws.Range("A78").Value = -19.533
$ws.Range("D81").Value = -7.342000000000001
$ws.Range("E86").Value = 16.294
$ws.Range("B87").Value = 5.140000000000001
$ws.Range("B88").Value = 5.948000000000001
$ws.Range("D92").Value = -7.154000000000001
$ws.Range("C95").Value = -11.663
$ws.Range("B96").Value = 6.368
$ws.Range("E96").Value = 16.666
$ws.Range("B98").Value = 5.556
$ws.Range("B101").Value = 8.659000000000001
$ws.Range("C101").Value = -12.643
$ws.Range("B102").Value = 7.739999999999999
